$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") rows 2 through 115 currently hold the date serial
# 45189 (2023-09-20). Bump every one of them forward by one day to 45190
# (2023-09-21).
for ($row = 2; $row -le 115; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45189) {
        $cell.Value2 = 45190
    }
}
